$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.09968638420105
$ws.Range("B1").Value = 2.798731565475464
$ws.Range("C1").Value = 5.108551025390625
$ws.Range("D1").Value = 2.095155954360962
$ws.Range("E1").Value = 1.169235706329346
